$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "UniformA-HW25.xpc" to "UniformA"
$ws.Name = "UniformA"

# Add new row 16, mirroring the formatting/pattern of row 15
# Copy formats of row 15 (A15:P15) down to row 16 first
$ws.Range("A15:P15").Copy()
$ws.Range("A16:P16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the values for the new row
$ws.Range("A16").Value = 14
$ws.Range("B16").Value2 = $ws.Range("B15").Value2

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 1
$ws.Range("P16").Value = 1
